# Apply updated 2025 (column L) crime totals for data through 2025-05-12.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2211
$ws.Range('L3').Value = 2228
$ws.Range('L4').Value = 613
$ws.Range('L5').Value = 132
$ws.Range('L6').Value = 2011
$ws.Range('L7').Value = 7195

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 57
$ws.Range('L5').Value = 24
$ws.Range('L7').Value = 233
$ws.Range('L8').Value = 456
$ws.Range('L9').Value = 46
$ws.Range('L11').Value = 129
$ws.Range('L12').Value = 18
$ws.Range('L14').Value = 38
$ws.Range('L20').Value = 185
$ws.Range('L29').Value = 367
$ws.Range('L33').Value = 323
$ws.Range('L34').Value = 48
$ws.Range('L37').Value = 260
$ws.Range('L42').Value = 223
$ws.Range('L45').Value = 13
$ws.Range('L47').Value = 54
$ws.Range('L51').Value = 83
$ws.Range('L53').Value = 91
$ws.Range('L54').Value = 146
$ws.Range('L55').Value = 65
$ws.Range('L63').Value = 25
$ws.Range('L64').Value = 50
$ws.Range('L65').Value = 139
$ws.Range('L67').Value = 256
$ws.Range('L68').Value = 21
$ws.Range('L73').Value = 55
$ws.Range('L76').Value = 82
$ws.Range('L77').Value = 43
$ws.Range('L79').Value = 198
$ws.Range('L80').Value = 26
$ws.Range('L82').Value = 17
$ws.Range('L83').Value = 176
$ws.Range('L85').Value = 379
$ws.Range('L86').Value = 57
$ws.Range('L91').Value = 102
$ws.Range('L94').Value = 85
$ws.Range('L96').Value = 69
$ws.Range('L101').Value = 7195

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L2').Value = 15
$ws.Range('L7').Value = 38

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 15
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 69

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L3').Value = 73
$ws.Range('L7').Value = 233

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 44
$ws.Range('L3').Value = 38
$ws.Range('L7').Value = 129

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 116
$ws.Range('L3').Value = 157
$ws.Range('L7').Value = 379

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 91

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 130
$ws.Range('L3').Value = 154
$ws.Range('L7').Value = 456

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 72
$ws.Range('L7').Value = 176

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 83
$ws.Range('L4').Value = 19
$ws.Range('L6').Value = 113
$ws.Range('L7').Value = 323

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 80
$ws.Range('L3').Value = 77
$ws.Range('L5').Value = 10
$ws.Range('L6').Value = 75
$ws.Range('L7').Value = 260

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 52
$ws.Range('L6').Value = 36
$ws.Range('L7').Value = 139

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 84
$ws.Range('L6').Value = 69
$ws.Range('L7').Value = 256

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L3').Value = 27
$ws.Range('L7').Value = 146

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L3').Value = 129
$ws.Range('L4').Value = 15
$ws.Range('L7').Value = 367

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L3').Value = 14
$ws.Range('L4').Value = 13
$ws.Range('L6').Value = 39
$ws.Range('L7').Value = 82

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 57
$ws.Range('L3').Value = 65
$ws.Range('L4').Value = 22
$ws.Range('L6').Value = 74
$ws.Range('L7').Value = 223

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 27
$ws.Range('L7').Value = 65

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 36
$ws.Range('L7').Value = 102

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 65
$ws.Range('L7').Value = 198

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 50

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 59
$ws.Range('L7').Value = 185

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L2').Value = 12
$ws.Range('L3').Value = 14
$ws.Range('L7').Value = 48

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L4').Value = 13
$ws.Range('L6').Value = 26
$ws.Range('L7').Value = 85

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L2').Value = 22
$ws.Range('L7').Value = 54

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('L2').Value = 12
$ws.Range('L6').Value = 12
$ws.Range('L7').Value = 46

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 55

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L6').Value = 10
$ws.Range('L7').Value = 24

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 34
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L3').Value = 26
$ws.Range('L7').Value = 83

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('L6').Value = 6
$ws.Range('L7').Value = 21

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('L2').Value = 2
$ws.Range('L7').Value = 17

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('L2').Value = 13
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('L3').Value = 7
$ws.Range('L7').Value = 13

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('L6').Value = 14
$ws.Range('L7').Value = 26

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('L2').Value = 4
$ws.Range('L7').Value = 18
